# Update pay-period-21 hours: reduce carry-over hours to 0 for rows 6, 8, 10
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 - Dharam Pal
$ws.Range("F6").Value = 80
$ws.Range("G6").Value = 40
$ws.Range("H6").Value = 40
$ws.Range("I6").Value = ""

# Row 8 - Raheel Shahzad
$ws.Range("F8").Value = 81.5
$ws.Range("G8").Value = 40.75
$ws.Range("H8").Value = 40.75
$ws.Range("I8").Value = ""

# Row 10 - Yulia McCoy
$ws.Range("F10").Value = 50.25
$ws.Range("H10").Value = 40.25
$ws.Range("I10").Value = ""
